$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the Paragraph object containing the (first) occurrence of
# the given text. (Range.Paragraphs.Item(1) is unreliable on sub-ranges in
# this host, so resolve via Document.Paragraphs + position comparison.)
# ---------------------------------------------------------------------------
function Find-ParagraphByText($doc, $text) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if (-not $found) {
        return $null
    }
    $s = $rng.Start
    $e = $rng.End
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Start -le $s -and $p.Range.End -ge $e) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (was after "Check the payment
#    methods is functional" under "Define Scope of Testing").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Remove the "Install/Uninstall Testing" bullet paragraph entirely.
# ---------------------------------------------------------------------------
$target = Find-ParagraphByText $d "Install/Uninstall Testing"
if ($target -ne $null) {
    $target.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Append the new "Determine Test Deliverables" sub-bullets after the
#    "Determine Test Deliverables" paragraph.
# ---------------------------------------------------------------------------
$anchor = Find-ParagraphByText $d "Determine Test Deliverables"

$items = @(
    @(1, "Test deliverables are provided before testing phase."),
    @(2, "Test plans document."),
    @(2, "Test cases documents"),
    @(2, "Test Design specifications."),
    @(1, "Test deliverables are provided during the testing"),
    @(2, "Test Scripts"),
    @(2, "Simulators."),
    @(2, "Test Data"),
    @(2, "Test Traceability Matrix"),
    @(2, "Error logs and execution logs."),
    @(1, "Test deliverables are provided after the testing cycles is over."),
    @(2, "Test Results/reports"),
    @(2, "Defect Report"),
    @(2, "Installation/ Test procedures guidelines"),
    @(2, "Release notes")
)

$prev = $anchor
$newParas = @()
foreach ($item in $items) {
    $level = $item[0]
    $text = $item[1]

    $prev.Range.InsertParagraphAfter()
    $newPara = $prev.Next()

    $newPara.Range.Text = $text
    $newPara.Range.ListFormat.ListLevelNumber = $level + 1

    $newParas += $newPara
    $prev = $newPara
}

# ---------------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark at the start of the new "Release notes"
#    paragraph (the last item added above).
# ---------------------------------------------------------------------------
$releasePara = $newParas[$newParas.Count - 1]
$bmStart = $releasePara.Range.Start
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
